$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B33").Value = "Stellar"
$ws.Range("C33").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"

# D and E columns: use quote-prefix to force text, avoiding Excel's numeric auto-detection,
# then clear the resulting quote-prefix style so no stray formatting is left behind.
$ws.Range("D2").Value = "'27.691.02"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.84%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.754.46"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -2.52%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.25%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'324.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -3.96%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.10%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4509"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.61%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3709"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -1.73%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'45.21"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  +0.23%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.07505"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.69%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'1.123"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -2.01%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = "'  -0.32%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'21.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.10%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'6.177"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.49%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'7.290"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -2.86%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'1.752.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -2.73%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.00001073"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -1.63%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'88.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  +8.21%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.06219"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -7.84%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Value = "'17.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -1.64%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'6.169"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.65%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.5297"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -3.66%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'27.709.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -1.81%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'  -1.93%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.321"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -4.21%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'20.60"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.55%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'153.34"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +0.97%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'2.359"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -0.15%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'1.950.52"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -2.78%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'127.90"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -3.88%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = "'  -1.25%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'0.09309"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -2.26%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "'5.734"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -1.69%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'3.637"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -9.75%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'12.66"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  +4.82%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.2174"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -5.56%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "'  -1.65%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.097"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -2.89%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.6478"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -1.90%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.06116"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.68%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'1.198"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -2.86%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'7.962"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -5.50%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'1.416"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -4.64%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.000"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.13%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'  -3.44%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'3.758"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -2.69%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.5946"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.63%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'126.11"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -3.46%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("E50").Value = "'  -2.57%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.06898"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -3.46%  "
$ws.Range("E51").Style = "Normal"
